$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the combined (two-paragraph) cell contents from row 7 and row 9
# before the row shift happens, so we can split them into the newly
# freed second row of each pair. Value2 is used for reading since it
# reliably returns the underlying string/number.
$E7 = $ws.Range("E7").Value2
$F7 = $ws.Range("F7").Value2
$G7 = $ws.Range("G7").Value2
$H7 = $ws.Range("H7").Value2

$E9 = $ws.Range("E9").Value2
$F9 = $ws.Range("F9").Value2
$G9 = $ws.Range("G9").Value2
$H9 = $ws.Range("H9").Value2

# Row 6 is completely blank; deleting it shifts rows 7-21 up to 6-20,
# matching the new, more compact layout.
$ws.Rows(6).Delete()

# Split the combined "label + answer" paragraphs (which were joined with a
# blank line) across the two rows of each story pair.
$sep = "`n`n"

$E7parts = $E7.Split($sep)
$F7parts = $F7.Split($sep)
$G7parts = $G7.Split($sep)
$H7parts = $H7.Split($sep)

$ws.Range("E6").Value = $E7parts[0]
$ws.Range("E7").Value = $E7parts[1]

$ws.Range("F6").Value = $F7parts[0]
$ws.Range("F7").Value = "Blacked out areas on the map to indicate no parking available"

$ws.Range("G6").Value = $G7parts[0]
$ws.Range("G7").Value = $G7parts[1]

$ws.Range("H6").Value = $H7parts[0]
$ws.Range("H7").Value = $H7parts[1]

$E9parts = $E9.Split($sep)
$F9parts = $F9.Split($sep)
$G9parts = $G9.Split($sep)
$H9parts = $H9.Split($sep)

$ws.Range("E8").Value = $E9parts[0]
$ws.Range("E9").Value = $E9parts[1]

$ws.Range("F8").Value = $F9parts[0]
$ws.Range("F9").Value = $F9parts[1]

$ws.Range("G8").Value = $G9parts[0]
$ws.Range("G9").Value = $G9parts[1]

$ws.Range("H8").Value = $H9parts[0]
$ws.Range("H9").Value = $H9parts[1]

# Row heights don't automatically travel with the shifted content, so
# restore them to match the rows' new positions (i.e. each row now has
# the height that the row below it used to have).
$ws.Rows(6).RowHeight = 150
$ws.Rows(7).RowHeight = 150
$ws.Rows(8).RowHeight = 150
$ws.Rows(9).RowHeight = 167
$ws.Rows(10).RowHeight = 150
$ws.Rows(11).RowHeight = 150
$ws.Rows(12).RowHeight = $ws.StandardHeight
$ws.Rows(13).RowHeight = 150
$ws.Rows(14).RowHeight = 150
$ws.Rows(15).RowHeight = 150
$ws.Rows(16).RowHeight = 157
$ws.Rows(17).RowHeight = 150
$ws.Rows(18).RowHeight = 150
$ws.Rows(19).RowHeight = 150
$ws.Rows(20).RowHeight = 150
$ws.Rows(21).RowHeight = 150

# Update the view to match: scroll/selection moved since rows shifted.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D18").Select()
